$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.264.09"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.230.07"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.91"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "78.39"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.14"
$ws.Range("E10").Value = "  +5.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0924"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.08"
$ws.Range("E12").Value = "  +3.00%  "
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.563.18"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.63"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.226.86"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.794"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.199.29"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.31"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.32"
$ws.Range("E22").Value = "  +6.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.13"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.32"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "42.50"
$ws.Range("E26").Value = "  +9.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.84"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.34"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.63"
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0872"
$ws.Range("E33").Value = "  +9.83%  "
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0369"
$ws.Range("E36").Value = "  +12.93%  "
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.45"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.25"
$ws.Range("E39").Value = "  +7.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("E40").Value = "  +19.84%  "
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.204"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "61.83"
$ws.Range("E43").Value = "  +3.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.34"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.84"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.58"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.478"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("E51").Value = "  +22.57%  "
